# Preset Location Values Update
#
# - Cell C2 on "Configuration" holds the notification e-mail address used
#   by the "Equals 3 Star and Daily" preset. Its address changes from
#   "chappel.mann+stl2@gmail.com" to "chappel.mann+stc@gmail.com", and the
#   mailto: hyperlink that used to be attached to that cell is removed.
#   The other rows (C3, C4, C5, C6, C7) keep their own mailto hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Configuration")

# --- Update the preset e-mail value in C2 ------------------------------
$ws.Range("C2").Value = "chappel.mann+stc@gmail.com"

# --- Remove the hyperlink that lived on C2 -----------------------------
# The host's Hyperlinks collection only supports clearing every hyperlink
# on the sheet at once (a per-item .Delete() is a no-op here), so record
# the mailto targets that must survive, wipe the collection, then recreate
# just those - leaving C2 without a hyperlink, as intended.
$survivors = @(
    @("C3", "mailto:sajith@gmail.com"),
    @("C6", "mailto:spillai@dacgroup.com"),
    @("C7", "mailto:spillai@dacgroup.com"),
    @("C4", "mailto:spillai@dacgroup.com"),
    @("C5", "mailto:spillai@dacgroup.com")
)

$ws.Hyperlinks.Delete()

foreach ($entry in $survivors) {
    $cell = $ws.Range($entry[0])
    $underlineBefore = $cell.Font.Underline
    $ws.Hyperlinks.Add($cell, $entry[1])
    # Recreating the hyperlink reapplies Excel's default (underlined)
    # hyperlink font; restore the cell's original underline state so the
    # existing formatting on these untouched rows is left unaffected.
    $cell.Font.Underline = $underlineBefore
}
